$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.5017037157228622, 0.5017037157228622, 0.5451339915373766, 0.148407217410269, 3.382615389267453, 6163)
    3 = @(0.6510359869138496, 0.6510359869138496, 0.5931445603576752, 0.06505041191066246, 3.321614489972583, 2751)
    4 = @(0.630192878338279, 0.630192878338279, 0.5702777544684065, 0.1345874138558461, 3.506261720912169, 5392)
    5 = @(0.3377982419422352, 0.3377982419422352, 0.3779859484777517, 0.07477551683647496, 4.396762830076731, 2389)
    6 = @(0.5311121999020089, 0.5311121999020089, 0.4716119208179247, 0.06218648568681087, 4.279990209406025, 2041)
    7 = @(0.4044054747647562, 0.4044054747647562, 0.4368718955758346, 0.1375679557111161, 4.132687312799809, 4676)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
